$wb = $excel.ActiveWorkbook

# --- "Login Page" sheet (sheet1.xml) ---
$ws1 = $wb.Worksheets.Item("Login Page")
$ws1.Activate()

# Change A2 value (text) from "9958592171" to "2537461015"
# Leading apostrophe keeps it stored as text (shared string), matching the
# cell's existing quote-prefixed text style instead of being read as a number.
$ws1.Range("A2").Value = "'2537461015"

# Column A width change (also drops the bestFit autofit flag)
$ws1.Columns("A").ColumnWidth = 18.75

# Update the view's active cell/selection
$ws1.Range("C7").Select()

# --- "Selected Vehicle Page" sheet (sheet2.xml) ---
$ws2 = $wb.Worksheets.Item("Selected Vehicle Page")
$ws2.Activate()
$ws2.Range("A11").Select()
